$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 81; this shifts the existing rows 81-114 down to 82-115,
# carrying their formatting (incl. the date-formatted D column) along with them.
$ws.Rows.Item(81).Insert()

# The row that used to be row 81 is now row 82 - duplicate its row of data into the
# freshly inserted row 81, then overwrite the date (column D) with the new value.
$src = $ws.Range("A82:R82")
$dst = $ws.Range("A81:R81")
$dst.Value = $src.Value()

$ws.Cells.Item(81, 4).Value = 44466
